# Updates cryptos list (Price / Volume(1h) columns, and two coin-row swaps)
# per the "Updated cryptos list" GitHub Actions commit.
# NumberFormat '@' (Text) is set on Price cells before assigning so that
# price strings like "1.00", "0.490", "96.330.37" keep their exact digits
# instead of being auto-coerced to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '96.330.37'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.569.51'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.89'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '656.02'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.54'
$ws.Range('E7').Value = '  +4.77%  '
$ws.Range('E8').Value = '  -1.78%  '
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.04'
$ws.Range('E10').Value = '  +2.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.567.63'
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.16'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.35'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.254.15'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.269.20'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.579.90'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.75'
$ws.Range('E19').Value = '  -2.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.54'
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.72'
$ws.Range('E21').Value = '  -2.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.490'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '510.46'
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.43'
$ws.Range('E24').Value = '  -2.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000199'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.80'
$ws.Range('E26').Value = '  +1.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '96.22'
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.74'
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.761.07'
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.98'
$ws.Range('E30').Value = '  -6.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.148'
$ws.Range('E31').Value = '  +6.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.46'
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.184'
$ws.Range('E34').Value = '  +3.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.54'
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.560'
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '603.00'
$ws.Range('E38').Value = '  +6.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.51'
$ws.Range('E39').Value = '  +2.47%  '
$ws.Range('E40').Value = '  +8.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.150'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.903'
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.82'
$ws.Range('E44').Value = '  +4.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.70'
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.50'
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '34.03'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0417'
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.62'
$ws.Range('E50').Value = '  +5.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.38'
$ws.Range('E51').Value = '  -1.97%  '
